$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U")

# --- New row 3: duplicate row 2's formatting + values first, while row 2
# still carries its original formatting (including U2's yellow/text
# style), so the copy is a faithful clone. ---------------------------
$ws.Range("A2:U2").Copy()
$ws.Range("A3:U3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

foreach ($c in $cols) {
    $src = $ws.Range($c + "2")
    $dst = $ws.Range($c + "3")
    $dst.Value = $src.Value2
}

# --- Row 2 edit: Cliente (A2): 3534375 -> 5802202 (force text so the
# numeric-looking string stays a shared string, same as before). ------
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "5802202"

# --- Row 3 field-specific overrides for the new proposal row (force
# text so these numeric-looking values are stored as shared strings,
# like the rest of the row). ------------------------------------------
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "19499545"

$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "20"

$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "20"

# --- Numero Propuesta updates: gggg -> 4873648 (row 2), new 4873649
# (row 3). Both drop the yellow highlight / text-format the cell
# previously had so they end up unstyled, matching the header cell U1. -
$ws.Range("U2").NumberFormat = "@"
$ws.Range("U2").Value = "4873648"
$ws.Range("U2").ClearFormats()

$ws.Range("U3").NumberFormat = "@"
$ws.Range("U3").Value = "4873649"
$ws.Range("U3").ClearFormats()

# --- Selection ---------------------------------------------------------
$ws.Range("H3").Select()
